$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update header row labels (row 1)
$ws.Range("B1").Value = "Government-Cadw"
$ws.Range("C1").Value = "Government-Local_Authority"
$ws.Range("D1").Value = "Government-National"
$ws.Range("E1").Value = "Government-Other"
$ws.Range("F1").Value = "Independent-English_Heritage"
$ws.Range("G1").Value = "Independent-Historic_Environment_Scotland"
$ws.Range("H1").Value = "Independent-National_Trust"
$ws.Range("I1").Value = "Independent-National_Trust_for_Scotland"
$ws.Range("J1").Value = "Independent-Not_for_profit"
$ws.Range("K1").Value = "Independent-Private"
$ws.Range("L1").Value = "Independent-Unknown"
$ws.Range("M1").Value = "University"
$ws.Range("N1").Value = "Unknown"

# Update row 2 (Accredited) values
$ws.Range("C2").Value = 583
$ws.Range("D2").Value = 65
$ws.Range("E2").Value = 1
$ws.Range("F2").Value = 33
$ws.Range("G2").Value = 12
$ws.Range("H2").Value = 143
$ws.Range("I2").Value = 11
$ws.Range("J2").Value = 799
$ws.Range("M2").Value = 73

# Update row 3 (Unaccredited) values
$ws.Range("C3").Value = 339
$ws.Range("D3").Value = 17
$ws.Range("E3").Value = 9
$ws.Range("F3").Value = 20
$ws.Range("G3").Value = 9
$ws.Range("H3").Value = 42
$ws.Range("I3").Value = 16
$ws.Range("J3").Value = 935
$ws.Range("K3").Value = 751
$ws.Range("L3").Value = 220
$ws.Range("M3").Value = 37
$ws.Range("N3").Value = 110
